# Add ability to load external files into TPA
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the sheet from "Sheet1" to "All"
$ws.Name = "All"

# Update the style used by K6/L6 (remove fill/border application, keep font+alignment)
$ws.Range("K6:L6").Borders.LineStyle = -4142
$ws.Range("K6:L6").Interior.Pattern = -4142

# Update the active selection on the sheet to A20 (below the data range)
$ws.Range("A20").Select()
